$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (interested-count) column F for several rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F11").Value = 116
$ws1.Range("F16").Value = 7294
$ws1.Range("F25").Value = 1730
$ws1.Range("F28").Value = 6134
$ws1.Range("F34").Value = 6329
$ws1.Range("F46").Value = 409
$ws1.Range("F47").Value = 2117

# Sheet "全部类型" - same column F updates (plus one row reset to 0)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 116
$ws4.Range("F16").Value = 7294
$ws4.Range("F23").Value = 1730
$ws4.Range("F28").Value = 6134
$ws4.Range("F35").Value = 6329
$ws4.Range("F45").Value = 0
$ws4.Range("F46").Value = 409
$ws4.Range("F48").Value = 2117
